$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Columns updated per row, in order: C, D, E, F, G, H, J, M, N
$colLetters = @("C","D","E","F","G","H","J","M","N")

$data = @(
    @(2, 4.914950792835458, 7.232443616329953, 11.01409948054581, 44.73207290059849, 59.60481127584346, 20.64422773698488, 9.078379738954178, 29.00560379347686, 17.10020804765148),
    @(3, 4.751718583703591, 7.110324789805595, 10.99251289653588, 44.25420883016235, 58.31999297991726, 20.51702392276673, 9.116171610652568, 28.19426810227755, 16.90396122823324),
    @(4, 4.650572766439653, 7.033423505892118, 10.97969372185768, 43.97693174143379, 57.54455822103629, 20.44600483941624, 9.140449827223883, 27.6879812878671, 16.78458359239772),
    @(5, 4.609195281905808, 7.001618968934259, 10.97458091249823, 43.86809600336566, 57.23242493420726, 20.41885875275551, 9.150614587140147, 27.47991698697832, 16.73626877775123),
    @(6, 4.602317122286293, 6.996310148957993, 10.97373869383729, 43.8502776865602, 57.18084271625633, 20.41445988586528, 9.152318850889248, 27.44527143103032, 16.72826771590188),
    @(7, 4.650015283736959, 7.032996446459363, 10.97962431663458, 43.97544698924636, 57.54033240955236, 20.44563145580851, 9.140585813150395, 27.68518193650544, 16.78393058797123),
    @(8, 4.858900944144376, 7.190744969949384, 11.00656716113006, 44.56402022134782, 59.15930364935883, 20.59890646996362, 9.09118820067801, 28.72771217784572, 17.03234075554067),
    @(9, 5.25832246952499, 7.484185637731549, 11.06277892269225, 45.84131821282696, 62.41972127004297, 20.95494830706694, 9.002786524679824, 30.69568103137626, 17.52609470345619),
    @(10, 5.542054347930895, 7.689239182003995, 11.10603520644482, 46.84761597670753, 64.8396326717583, 21.24913846262579, 8.942923259236585, 32.0803250166916, 17.88989018816469),
    @(11, 5.668388917695323, 7.780072405592985, 11.12611620935524, 47.31847219309314, 65.94023958133796, 21.38973443951414, 8.916777651607049, 32.69425462796974, 18.05498905860091),
    @(12, 5.715789680896156, 7.814104042958537, 11.13377634346062, 47.49850797267273, 66.35654023365794, 21.44391678936293, 8.907031947398233, 32.92424298858464, 18.11740004507611),
    @(13, 5.705601297033983, 7.806791151311197, 11.13212414500853, 47.45965938185812, 66.26691037466608, 21.43220629208573, 8.90912398101613, 32.87482474870473, 18.10396435269273),
    @(14, 5.67229765778489, 7.78287959477825, 11.12674530535245, 47.33324998113607, 65.97450086393349, 21.39417336882683, 8.915972765530586, 32.7132269555328, 18.06012611186794),
    @(15, 5.651839720494481, 7.768185176190119, 11.12345782585229, 47.25604175931329, 65.79531695951111, 21.37099877189639, 8.920188003754085, 32.61391332819531, 18.03325829585248),
    @(16, 5.533738980398029, 7.683252644828338, 11.10473080424206, 46.81709507028045, 64.76766249676483, 21.240083726376, 8.944653699625032, 32.03986477061195, 17.87908821389704),
    @(17, 5.460552818440663, 7.630513050949779, 11.09334422974061, 46.55106213276255, 64.13684975627501, 21.16148323970889, 8.959940060358335, 31.68347065664301, 17.78437182555642),
    @(18, 5.41820290930397, 7.599949071430931, 11.08683308322111, 46.39928792677541, 63.77402704897352, 21.11691333494544, 8.968834706660617, 31.47699138543278, 17.72985944756974),
    @(19, 5.403821616448378, 7.589561618988784, 11.08463513042381, 46.34811740734627, 63.65119611566303, 21.10193340298723, 8.971863892759711, 31.40683140845222, 17.71139824942066),
    @(20, 5.468370367709674, 7.636151126907242, 11.09455242107288, 46.57925443515051, 64.20400371310215, 21.16978448056546, 8.958302218217121, 31.72156521890467, 17.7944584466223),
    @(21, 5.682092009332891, 7.789913003062837, 11.12832370118877, 47.37033366067705, 66.06040493570984, 21.4053192641774, 8.91395691226994, 32.76076131050127, 18.07300581316369),
    @(22, 5.819192596156414, 7.888272679346156, 11.15071995959512, 47.89737880357075, 67.27071220499988, 21.56472776413016, 8.885878011806247, 33.42532322343214, 18.25439789848481),
    @(23, 5.746269246439117, 7.835975516020695, 11.13873765308359, 47.61521706503481, 66.62515911736146, 21.47915861747452, 8.900781982533418, 33.0720320440926, 18.1576621913069),
    @(24, 5.464836904249327, 7.63360290954037, 11.09400608813628, 46.56650503446433, 64.17364387386601, 21.16602956053341, 8.959042355417008, 31.70434758434971, 17.78989846402395),
    @(25, 5.151743491392846, 7.406594575505938, 11.04722085599353, 45.48330728333806, 61.53143361775621, 20.85280392962137, 9.025802939837128, 30.17306492948946, 17.39212408504709)
)

foreach ($entry in $data) {
    $row = $entry[0]
    for ($i = 0; $i -lt $colLetters.Length; $i++) {
        $col = $colLetters[$i]
        $value = $entry[$i + 1]
        $ws.Range("$col$row").Value = $value
    }
}
